# Add two new batches (30 and 31) of devices to the Master Device table:
# 5 rows each for Finger Print Scanner / IRIS Scanner / Web Camera /
# Document Scanner / Printer, mirroring the existing rows' layout/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: id
$ws.Cells.Item(147, 1).Value = 3000166
$ws.Cells.Item(148, 1).Value = 3000167
$ws.Cells.Item(149, 1).Value = 3000168
$ws.Cells.Item(150, 1).Value = 3000169
$ws.Cells.Item(151, 1).Value = 3000170
$ws.Cells.Item(152, 1).Value = 3000171
$ws.Cells.Item(153, 1).Value = 3000172
$ws.Cells.Item(154, 1).Value = 3000173
$ws.Cells.Item(155, 1).Value = 3000174
$ws.Cells.Item(156, 1).Value = 3000175

# Column B: name (batch 30)
$ws.Cells.Item(147, 2).Value = "Finger Print Scanner 30"
$ws.Cells.Item(148, 2).Value = "IRIS Scanner 30"
$ws.Cells.Item(149, 2).Value = "Web Camera 30"
$ws.Cells.Item(150, 2).Value = "Document Scanner 30"
$ws.Cells.Item(151, 2).Value = "Printer 30"

# Column C: mac_address (batch 30)
$ws.Cells.Item(147, 3).Value = "D6-15-AC-80-6B-86"
$ws.Cells.Item(148, 3).Value = "6D-58-E2-DF-74-34"
$ws.Cells.Item(149, 3).Value = "E2-A8-56-86-15-30"
$ws.Cells.Item(150, 3).Value = "72-E8-B9-FD-63-65"
$ws.Cells.Item(151, 3).Value = "D3-F3-A4-50-AD-12"

# Column D: serial_num (batch 30)
$ws.Cells.Item(147, 4).Value = "BS563Q2230814"
$ws.Cells.Item(148, 4).Value = "BS563Q2230815"
$ws.Cells.Item(149, 4).Value = "BS563Q2230816"
$ws.Cells.Item(150, 4).Value = "BS563Q2230817"
$ws.Cells.Item(151, 4).Value = "BS563Q2230818"

# Column B: name (batch 31)
$ws.Cells.Item(152, 2).Value = "Finger Print Scanner 31"
$ws.Cells.Item(153, 2).Value = "IRIS Scanner 31"
$ws.Cells.Item(154, 2).Value = "Web Camera 31"
$ws.Cells.Item(155, 2).Value = "Document Scanner 31"
$ws.Cells.Item(156, 2).Value = "Printer 31"

# Column D: serial_num (batch 31)
$ws.Cells.Item(152, 4).Value = "BS563Q2230819"
$ws.Cells.Item(153, 4).Value = "BS563Q2230820"
$ws.Cells.Item(154, 4).Value = "BS563Q2230821"
$ws.Cells.Item(155, 4).Value = "BS563Q2230822"
$ws.Cells.Item(156, 4).Value = "BS563Q2230823"

# Column C: mac_address (batch 31)
$ws.Cells.Item(152, 3).Value = "06-16-D0-0B-A6-E4"
$ws.Cells.Item(153, 3).Value = "21-78-45-AC-E9-20"
$ws.Cells.Item(154, 3).Value = "3C-E8-87-99-DB-FA"
$ws.Cells.Item(155, 3).Value = "BF-55-53-98-40-08"
$ws.Cells.Item(156, 3).Value = "5A-43-36-46-22-EB"

# Column F: dspec_id
$ws.Cells.Item(147, 6).Value = 165
$ws.Cells.Item(148, 6).Value = 327
$ws.Cells.Item(149, 6).Value = 736
$ws.Cells.Item(150, 6).Value = 801
$ws.Cells.Item(151, 6).Value = 920
$ws.Cells.Item(152, 6).Value = 165
$ws.Cells.Item(153, 6).Value = 327
$ws.Cells.Item(154, 6).Value = 736
$ws.Cells.Item(155, 6).Value = 801
$ws.Cells.Item(156, 6).Value = 920

# Columns G, H, I, J: constant values (lang_code, is_active, cr_by, cr_dtimes)
$ws.Cells.Item(147, 7).Value = "eng"
$ws.Cells.Item(147, 8).Value = $true
$ws.Cells.Item(147, 8).HorizontalAlignment = -4131
$ws.Cells.Item(147, 9).Value = "superadmin"
$ws.Cells.Item(147, 10).Value = "now()"
$ws.Cells.Item(148, 7).Value = "eng"
$ws.Cells.Item(148, 8).Value = $true
$ws.Cells.Item(148, 8).HorizontalAlignment = -4131
$ws.Cells.Item(148, 9).Value = "superadmin"
$ws.Cells.Item(148, 10).Value = "now()"
$ws.Cells.Item(149, 7).Value = "eng"
$ws.Cells.Item(149, 8).Value = $true
$ws.Cells.Item(149, 8).HorizontalAlignment = -4131
$ws.Cells.Item(149, 9).Value = "superadmin"
$ws.Cells.Item(149, 10).Value = "now()"
$ws.Cells.Item(150, 7).Value = "eng"
$ws.Cells.Item(150, 8).Value = $true
$ws.Cells.Item(150, 8).HorizontalAlignment = -4131
$ws.Cells.Item(150, 9).Value = "superadmin"
$ws.Cells.Item(150, 10).Value = "now()"
$ws.Cells.Item(151, 7).Value = "eng"
$ws.Cells.Item(151, 8).Value = $true
$ws.Cells.Item(151, 8).HorizontalAlignment = -4131
$ws.Cells.Item(151, 9).Value = "superadmin"
$ws.Cells.Item(151, 10).Value = "now()"
$ws.Cells.Item(152, 7).Value = "eng"
$ws.Cells.Item(152, 8).Value = $true
$ws.Cells.Item(152, 8).HorizontalAlignment = -4131
$ws.Cells.Item(152, 9).Value = "superadmin"
$ws.Cells.Item(152, 10).Value = "now()"
$ws.Cells.Item(153, 7).Value = "eng"
$ws.Cells.Item(153, 8).Value = $true
$ws.Cells.Item(153, 8).HorizontalAlignment = -4131
$ws.Cells.Item(153, 9).Value = "superadmin"
$ws.Cells.Item(153, 10).Value = "now()"
$ws.Cells.Item(154, 7).Value = "eng"
$ws.Cells.Item(154, 8).Value = $true
$ws.Cells.Item(154, 8).HorizontalAlignment = -4131
$ws.Cells.Item(154, 9).Value = "superadmin"
$ws.Cells.Item(154, 10).Value = "now()"
$ws.Cells.Item(155, 7).Value = "eng"
$ws.Cells.Item(155, 8).Value = $true
$ws.Cells.Item(155, 8).HorizontalAlignment = -4131
$ws.Cells.Item(155, 9).Value = "superadmin"
$ws.Cells.Item(155, 10).Value = "now()"
$ws.Cells.Item(156, 7).Value = "eng"
$ws.Cells.Item(156, 8).Value = $true
$ws.Cells.Item(156, 8).HorizontalAlignment = -4131
$ws.Cells.Item(156, 9).Value = "superadmin"
$ws.Cells.Item(156, 10).Value = "now()"

# Update selection / scroll position to match the final view state
[void]$ws.Range("E156").Select()
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 1

# Bump print quality (mirrors pageSetup verticalDpi going from 0 to 300)
$ws.PageSetup.PrintQuality(2) = 300
